# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue 'D2' '66.767.62'
Set-TextValue 'E2' '  -0.25%  '
Set-TextValue 'D3' '3.067.25'
Set-TextValue 'E3' '  -1.58%  '
Set-TextValue 'E4' '  +0.06%  '
Set-TextValue 'D5' '575.87'
Set-TextValue 'E5' '  -0.45%  '
Set-TextValue 'D6' '168.07'
Set-TextValue 'E6' '  -2.46%  '
Set-TextValue 'D8' '3.062.36'
Set-TextValue 'E8' '  -1.63%  '
Set-TextValue 'D9' '0.511'
Set-TextValue 'E9' '  -2.09%  '
Set-TextValue 'D10' '6.38'
Set-TextValue 'E10' '  -1.25%  '
Set-TextValue 'E11' '  -2.25%  '
Set-TextValue 'D12' '0.469'
Set-TextValue 'E12' '  -2.86%  '
Set-TextValue 'D13' '0.0000240'
Set-TextValue 'E13' '  -2.57%  '
Set-TextValue 'D14' '35.72'
Set-TextValue 'E14' '  -4.05%  '
Set-TextValue 'D16' '66.784.83'
Set-TextValue 'E16' '  -0.10%  '
Set-TextValue 'D17' '3.579.76'
Set-TextValue 'E17' '  -1.42%  '
Set-TextValue 'D18' '6.99'
Set-TextValue 'E18' '  -2.14%  '
Set-TextValue 'D19' '16.83'
Set-TextValue 'E19' '  +2.66%  '
Set-TextValue 'D20' '3.072.69'
Set-TextValue 'E20' '  -1.33%  '
Set-TextValue 'D21' '488.69'
Set-TextValue 'E21' '  +2.47%  '
Set-TextValue 'D22' '0.688'
Set-TextValue 'E22' '  -3.69%  '
Set-TextValue 'D23' '7.68'
Set-TextValue 'E23' '  -4.20%  '
Set-TextValue 'D24' '82.79'
Set-TextValue 'E24' '  -1.51%  '
Set-TextValue 'D25' '12.69'
Set-TextValue 'E25' '  -6.32%  '
Set-TextValue 'D26' '2.21'
Set-TextValue 'E26' '  -4.18%  '
Set-TextValue 'D27' '10.19'
Set-TextValue 'E27' '  +1.83%  '
Set-TextValue 'E28' '  +0.02%  '
Set-TextValue 'D29' '7.78'
Set-TextValue 'E29' '  -1.53%  '
Set-TextValue 'D30' '2.27'
Set-TextValue 'E30' '  -5.52%  '
Set-TextValue 'D31' '2.61'
Set-TextValue 'E31' '  -2.37%  '
Set-TextValue 'D32' '27.58'
Set-TextValue 'E33' '  -3.20%  '
Set-TextValue 'D34' '0.0₃0910'
Set-TextValue 'E34' '  -3.58%  '
Set-TextValue 'E35' '  +0.10%  '
Set-TextValue 'D36' '5.62'
Set-TextValue 'E36' '  -4.37%  '
Set-TextValue 'D37' '0.948'
Set-TextValue 'E37' '  -3.10%  '
Set-TextValue 'D38' '46.79'
Set-TextValue 'E38' '  -1.02%  '
Set-TextValue 'D39' '0.122'
Set-TextValue 'E39' '  +0.48%  '
Set-TextValue 'D40' '1.97'
Set-TextValue 'E40' '  -5.37%  '
Set-TextValue 'D41' '0.300'
Set-TextValue 'E41' '  -3.33%  '
Set-TextValue 'D42' '8.32'
Set-TextValue 'E42' '  -4.19%  '
Set-TextValue 'D43' '2.756.06'
Set-TextValue 'E43' '  -2.03%  '
Set-TextValue 'D44' '370.59'
Set-TextValue 'E44' '  -2.75%  '
Set-TextValue 'D45' '0.0345'
Set-TextValue 'E45' '  -3.07%  '
Set-TextValue 'D46' '135.91'
Set-TextValue 'D47' '2.48'
Set-TextValue 'E47' '  -4.58%  '
Set-TextValue 'D49' '24.39'
Set-TextValue 'E49' '  -1.80%  '
Set-TextValue 'E50' '  -2.09%  '
Set-TextValue 'E51' '  -1.88%  '
